# Applies the updated crypto price/volume snapshot described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'43.103.74"
$ws.Range("E2").Value = "'  +0.27%  "

# Row 3
$ws.Range("D3").Value = "'2.219.98"
$ws.Range("E3").Value = "'  -0.68%  "

# Row 4
$ws.Range("E4").Value = "'  -0.05%  "

# Row 5
$ws.Range("D5").Value = "'257.92"
$ws.Range("E5").Value = "'  +2.39%  "

# Row 6
$ws.Range("D6").Value = "'0.626"
$ws.Range("E6").Value = "'  +2.07%  "

# Row 7
$ws.Range("D7").Value = "'77.56"
$ws.Range("E7").Value = "'  +3.32%  "

# Row 8
$ws.Range("E8").Value = "'  -0.07%  "

# Row 9
$ws.Range("D9").Value = "'0.597"
$ws.Range("E9").Value = "'  +0.18%  "

# Row 10
$ws.Range("D10").Value = "'42.80"
$ws.Range("E10").Value = "'  +4.15%  "

# Row 11
$ws.Range("D11").Value = "'0.0916"
$ws.Range("E11").Value = "'  -0.81%  "

# Row 12
$ws.Range("D12").Value = "'7.03"
$ws.Range("E12").Value = "'  +2.24%  "

# Row 13
$ws.Range("E13").Value = "'  +1.45%  "

# Row 14
$ws.Range("D14").Value = "'2.543.78"
$ws.Range("E14").Value = "'  -0.90%  "

# Row 15
$ws.Range("D15").Value = "'14.52"
$ws.Range("E15").Value = "'  +0.07%  "

# Row 16
$ws.Range("D16").Value = "'2.225.48"
$ws.Range("E16").Value = "'  -0.58%  "

# Row 17
$ws.Range("D17").Value = "'0.789"
$ws.Range("E17").Value = "'  +0.10%  "

# Row 18
$ws.Range("D18").Value = "'43.015.66"
$ws.Range("E18").Value = "'  +0.29%  "

# Row 19
$ws.Range("E19").Value = "'  -0.36%  "

# Row 20
$ws.Range("D20").Value = "'71.22"
$ws.Range("E20").Value = "'  +0.03%  "

# Row 21
$ws.Range("D21").Value = "'5.99"
$ws.Range("E21").Value = "'  +0.29%  "

# Row 22
$ws.Range("E22").Value = "'  +7.72%  "

# Row 23
$ws.Range("D23").Value = "'230.39"
$ws.Range("E23").Value = "'  +0.37%  "

# Row 24
$ws.Range("D24").Value = "'9.20"
$ws.Range("E24").Value = "'  -1.63%  "

# Row 25
$ws.Range("E25").Value = "'  -0.11%  "

# Row 26
$ws.Range("D26").Value = "'42.70"
$ws.Range("E26").Value = "'  +8.28%  "

# Row 27
$ws.Range("D27").Value = "'10.80"
$ws.Range("E27").Value = "'  +1.18%  "

# Row 28
$ws.Range("D28").Value = "'3.33"
$ws.Range("E28").Value = "'  -2.18%  "

# Row 29
$ws.Range("D29").Value = "'2.21"
$ws.Range("E29").Value = "'  -0.47%  "

# Row 30
$ws.Range("E30").Value = "'  +2.13%  "

# Row 31
$ws.Range("D31").Value = "'173.23"
$ws.Range("E31").Value = "'  +1.18%  "

# Row 32
$ws.Range("D32").Value = "'20.41"
$ws.Range("E32").Value = "'  +1.11%  "

# Row 33
$ws.Range("D33").Value = "'0.0866"
$ws.Range("E33").Value = "'  +9.18%  "

# Row 34
$ws.Range("D34").Value = "'5.25"
$ws.Range("E34").Value = "'  +0.49%  "

# Row 35
$ws.Range("D35").Value = "'0.122"
$ws.Range("E35").Value = "'  +0.95%  "

# Row 36
$ws.Range("D36").Value = "'0.0364"
$ws.Range("E36").Value = "'  +11.56%  "

# Row 37
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").Value = "'4.48"
$ws.Range("E37").Value = "'  +0.40%  "

# Row 38
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").Value = "'0.107"
$ws.Range("E38").Value = "'  -4.31%  "

# Row 39
$ws.Range("D39").Value = "'13.05"
$ws.Range("E39").Value = "'  +4.82%  "

# Row 40
$ws.Range("D40").Value = "'2.94"
$ws.Range("E40").Value = "'  +21.51%  "

# Row 41
$ws.Range("D41").Value = "'2.12"
$ws.Range("E41").Value = "'  +1.27%  "

# Row 42
$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").Value = "'0.203"
$ws.Range("E42").Value = "'  -0.20%  "

# Row 43
$ws.Range("B43").Value = "MultiversX"
$ws.Range("C43").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D43").Value = "'61.32"
$ws.Range("E43").Value = "'  +2.96%  "

# Row 44
$ws.Range("D44").Value = "'5.29"
$ws.Range("E44").Value = "'  -1.35%  "

# Row 45
$ws.Range("D45").Value = "'103.38"
$ws.Range("E45").Value = "'  +0.07%  "

# Row 46
$ws.Range("D46").Value = "'8.46"
$ws.Range("E46").Value = "'  -2.06%  "

# Row 47
$ws.Range("D47").Value = "'0.472"
$ws.Range("E47").Value = "'  -2.62%  "

# Row 48
$ws.Range("D48").Value = "'0.0976"
$ws.Range("E48").Value = "'  -0.96%  "

# Row 49
$ws.Range("E49").Value = "'  +0.43%  "

# Row 50
$ws.Range("D50").Value = "'1.15"
$ws.Range("E50").Value = "'  +0.54%  "

# Row 51
$ws.Range("D51").Value = "'1.47"
$ws.Range("E51").Value = "'  +22.86%  "
